# Updated main GSC export data: the oldest day (2025-11-15) has aged out of
# the export window, so drop its row from the "Chart" data sheet. Excel
# shifts every subsequent row up by one, which is exactly what the source
# export does when it rolls the date window forward by a day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest date (2025-11-15). Deleting it shifts every later
# row (dates + Non-HTTPS/HTTPS counts) up by one, so the sheet now runs
# from 2025-11-16 through 2026-02-13 with 90 data rows instead of 91.
$ws.Rows.Item(2).Delete()
